# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
# Commit: Updated cryptos list on Sat Oct 28 14:42:48 UTC 2023 with GitHub Actions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.127.83"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.788.24"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.72"
$ws.Range("E5").Value = "  +0.31%  "

# Row 6
$ws.Range("E6").Value = "  -0.60%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.89"
$ws.Range("E8").Value = "  -1.02%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +1.14%  "

# Row 10
$ws.Range("E10").Value = "  -2.09%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("E11").Value = "  +0.87%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.046.92"
$ws.Range("E12").Value = "  +0.59%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.18"
$ws.Range("E13").Value = "  +2.36%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.788.64"
$ws.Range("E14").Value = "  -0.12%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.057.18"
$ws.Range("E15").Value = "  +0.38%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.620"
$ws.Range("E16").Value = "  +0.09%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  +1.07%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.12"
$ws.Range("E18").Value = "  +0.40%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.81"
$ws.Range("E19").Value = "  +1.25%  "

# Row 20
$ws.Range("E20").Value = "  -0.64%  "

# Row 21
$ws.Range("E21").Value = "  -0.03%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.86"
$ws.Range("E22").Value = "  +1.88%  "

# Row 23
$ws.Range("E23").Value = "  +0.84%  "

# Row 24
$ws.Range("E24").Value = "  -1.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.61"
$ws.Range("E25").Value = "  +0.89%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  +1.23%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("E28").Value = "  +0.68%  "

# Row 29
$ws.Range("E29").Value = "  +0.12%  "

# Row 30
$ws.Range("E30").Value = "  -0.83%  "

# Row 31
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("E32").Value = "  +0.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.60"
$ws.Range("E33").Value = "  +3.10%  "

# Row 34
$ws.Range("E34").Value = "  -0.30%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.462.26"
$ws.Range("E35").Value = "  +5.04%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.43"
$ws.Range("E36").Value = "  +10.49%  "

# Row 37
$ws.Range("E37").Value = "  -1.20%  "

# Row 38
$ws.Range("E38").Value = "  +2.65%  "

# Row 39
$ws.Range("E39").Value = "  -0.17%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.30"
$ws.Range("E40").Value = "  +3.66%  "

# Row 41
$ws.Range("E41").Value = "  +0.31%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.919"
$ws.Range("E42").Value = "  +1.11%  "

# Row 43
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.46"
$ws.Range("E44").Value = "  +2.42%  "

# Row 45
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0509"
$ws.Range("E45").Value = "  +2.56%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.05"
$ws.Range("E46").Value = "  +4.21%  "

# Row 47
$ws.Range("E47").Value = "  -0.38%  "

# Row 48
$ws.Range("E48").Value = "  +0.51%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.948.26"
$ws.Range("E49").Value = "  +0.72%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.21"
$ws.Range("E50").Value = "  -1.50%  "

# Row 51
$ws.Range("E51").Value = "  -0.03%  "
